$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D ("Status"), shifting old D..H to E..I.
# This preserves the header row style (bold/border/center) on the
# shifted-right cells, including the newly exposed I1.
$ws.Columns.Item(4).Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "ISIN"
$ws.Cells.Item(1, 2).Value = "Stock Name"
$ws.Cells.Item(1, 3).Value = "Mutual Fund"
$ws.Cells.Item(1, 4).Value = "Status"
$ws.Cells.Item(1, 5).Value = "Jan_2026"
$ws.Cells.Item(1, 6).Value = "Dec_2025"
$ws.Cells.Item(1, 7).Value = "Oct_2025"
$ws.Cells.Item(1, 8).Value = "MoM"
$ws.Cells.Item(1, 9).Value = "QoQ"

# --- Data rows (rows 2-26) ---
# row 2
$ws.Cells.Item(2, 1).Value = "INE040A01034"
$ws.Cells.Item(2, 2).Value = "HDFC Bank Limited"
$ws.Cells.Item(2, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(2, 4).Value = "Fresh Entry"
$ws.Cells.Item(2, 5).Value = 9.187924
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 9.187924
$ws.Cells.Item(2, 9).Value = 9.187924

# row 3
$ws.Cells.Item(3, 1).Value = "INE202B01038"
$ws.Cells.Item(3, 2).Value = "Piramal Finance Ltd"
$ws.Cells.Item(3, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(3, 4).Value = "Reducing Consistently"
$ws.Cells.Item(3, 5).Value = 8.361349
$ws.Cells.Item(3, 6).Value = 10.172964
$ws.Cells.Item(3, 7).Value = 8.45973
$ws.Cells.Item(3, 8).Value = -1.811615
$ws.Cells.Item(3, 9).Value = -0.09838099999999983

# row 4
$ws.Cells.Item(4, 1).Value = "INE018A01030"
$ws.Cells.Item(4, 2).Value = "Larsen & Toubro Limited"
$ws.Cells.Item(4, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(4, 4).Value = "Adding Consistently"
$ws.Cells.Item(4, 5).Value = 7.955206
$ws.Cells.Item(4, 6).Value = 7.788444
$ws.Cells.Item(4, 7).Value = 7.500485
$ws.Cells.Item(4, 8).Value = 0.1667620000000003
$ws.Cells.Item(4, 9).Value = 0.4547210000000002

# row 5
$ws.Cells.Item(5, 1).Value = "INE423A01024"
$ws.Cells.Item(5, 2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(5, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(5, 4).Value = "Reducing Consistently"
$ws.Cells.Item(5, 5).Value = 7.508014
$ws.Cells.Item(5, 6).Value = 7.84677
$ws.Cells.Item(5, 7).Value = 8.601606
$ws.Cells.Item(5, 8).Value = -0.3387560000000001
$ws.Cells.Item(5, 9).Value = -1.093592

# row 6
$ws.Cells.Item(6, 1).Value = "INE090A01021"
$ws.Cells.Item(6, 2).Value = "ICICI Bank Limited"
$ws.Cells.Item(6, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(6, 4).Value = "Fresh Entry"
$ws.Cells.Item(6, 5).Value = 6.406086
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 6.406086
$ws.Cells.Item(6, 9).Value = 6.406086

# row 7
$ws.Cells.Item(7, 1).Value = "INE795G01014"
$ws.Cells.Item(7, 2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(7, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(7, 4).Value = "Adding Consistently"
$ws.Cells.Item(7, 5).Value = 5.947946
$ws.Cells.Item(7, 6).Value = 5.752245
$ws.Cells.Item(7, 7).Value = 5.477145
$ws.Cells.Item(7, 8).Value = 0.1957009999999997
$ws.Cells.Item(7, 9).Value = 0.4708009999999998

# row 8
$ws.Cells.Item(8, 1).Value = "INE364U01010"
$ws.Cells.Item(8, 2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(8, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(8, 4).Value = "Reducing"
$ws.Cells.Item(8, 5).Value = 5.256688
$ws.Cells.Item(8, 6).Value = 5.902581
$ws.Cells.Item(8, 7).Value = 4.278965
$ws.Cells.Item(8, 8).Value = -0.645893
$ws.Cells.Item(8, 9).Value = 0.9777229999999992

# row 9
$ws.Cells.Item(9, 1).Value = "INE406A01037"
$ws.Cells.Item(9, 2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(9, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(9, 4).Value = "Adding Consistently"
$ws.Cells.Item(9, 5).Value = 3.961293
$ws.Cells.Item(9, 6).Value = 3.658272
$ws.Cells.Item(9, 7).Value = 3.435944
$ws.Cells.Item(9, 8).Value = 0.3030209999999998
$ws.Cells.Item(9, 9).Value = 0.5253489999999998

# row 10
$ws.Cells.Item(10, 1).Value = "INE917I01010"
$ws.Cells.Item(10, 2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(10, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(10, 4).Value = "Adding Consistently"
$ws.Cells.Item(10, 5).Value = 3.561141
$ws.Cells.Item(10, 6).Value = 3.268366
$ws.Cells.Item(10, 7).Value = 3.034851
$ws.Cells.Item(10, 8).Value = 0.2927750000000002
$ws.Cells.Item(10, 9).Value = 0.5262899999999999

# row 11
$ws.Cells.Item(11, 1).Value = "INE237A01036"
$ws.Cells.Item(11, 2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(11, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(11, 4).Value = "Fresh Entry"
$ws.Cells.Item(11, 5).Value = 3.097165
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 3.097165
$ws.Cells.Item(11, 9).Value = 3.097165

# row 12
$ws.Cells.Item(12, 1).Value = "INE814H01029"
$ws.Cells.Item(12, 2).Value = "Adani Power Limited"
$ws.Cells.Item(12, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(12, 4).Value = "Reducing Consistently"
$ws.Cells.Item(12, 5).Value = 2.788864
$ws.Cells.Item(12, 6).Value = 4.034234
$ws.Cells.Item(12, 7).Value = 4.344794
$ws.Cells.Item(12, 8).Value = -1.24537
$ws.Cells.Item(12, 9).Value = -1.55593

# row 13
$ws.Cells.Item(13, 1).Value = "INE726G01019"
$ws.Cells.Item(13, 2).Value = "ICICI Prudential Life Insurance Co Ltd"
$ws.Cells.Item(13, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(13, 4).Value = "Fresh Entry"
$ws.Cells.Item(13, 5).Value = 2.351255
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 2.351255
$ws.Cells.Item(13, 9).Value = 2.351255

# row 14
$ws.Cells.Item(14, 1).Value = "INE200M01039"
$ws.Cells.Item(14, 2).Value = "Varun Beverages Limited"
$ws.Cells.Item(14, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(14, 4).Value = "Fresh Entry"
$ws.Cells.Item(14, 5).Value = 1.93914
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 1.93914
$ws.Cells.Item(14, 9).Value = 1.93914

# row 15
$ws.Cells.Item(15, 1).Value = "INE259A01022"
$ws.Cells.Item(15, 2).Value = "Colgate-Palmolive (India) Ltd"
$ws.Cells.Item(15, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(15, 4).Value = "Fresh Entry"
$ws.Cells.Item(15, 5).Value = 1.572049
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 1.572049
$ws.Cells.Item(15, 9).Value = 1.572049

# row 16
$ws.Cells.Item(16, 1).Value = "INE931S01010"
$ws.Cells.Item(16, 2).Value = "Adani Energy Solutions Limited"
$ws.Cells.Item(16, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(16, 4).Value = "Fresh Entry"
$ws.Cells.Item(16, 5).Value = 0.423678
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0.423678
$ws.Cells.Item(16, 9).Value = 0.423678

# row 17
$ws.Cells.Item(17, 1).Value = "INE424H01027"
$ws.Cells.Item(17, 2).Value = "SUN TV Network Limited"
$ws.Cells.Item(17, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(17, 4).Value = "Complete Exit"
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0.961596
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = -0.961596

# row 18
$ws.Cells.Item(18, 1).Value = "INE237A01028"
$ws.Cells.Item(18, 2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(18, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(18, 4).Value = "Complete Exit"
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 3.150556
$ws.Cells.Item(18, 7).Value = 2.935558
$ws.Cells.Item(18, 8).Value = -3.150556
$ws.Cells.Item(18, 9).Value = -2.935558

# row 19
$ws.Cells.Item(19, 1).Value = "INE019A01038"
$ws.Cells.Item(19, 2).Value = "JSW Steel Limited"
$ws.Cells.Item(19, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(19, 4).Value = "Complete Exit"
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 3.07038
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = -3.07038

# row 20
$ws.Cells.Item(20, 1).Value = "INE758E01017"
$ws.Cells.Item(20, 2).Value = "Jio Financial Services Limited"
$ws.Cells.Item(20, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(20, 4).Value = "Complete Exit"
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 8.168625
$ws.Cells.Item(20, 7).Value = 8.289438
$ws.Cells.Item(20, 8).Value = -8.168625
$ws.Cells.Item(20, 9).Value = -8.289438

# row 21
$ws.Cells.Item(21, 1).Value = "INE775A01035"
$ws.Cells.Item(21, 2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(21, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(21, 4).Value = "Complete Exit"
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 6.446333
$ws.Cells.Item(21, 7).Value = 2.869898
$ws.Cells.Item(21, 8).Value = -6.446333
$ws.Cells.Item(21, 9).Value = -2.869898

# row 22
$ws.Cells.Item(22, 1).Value = "INE326A01037"
$ws.Cells.Item(22, 2).Value = "Lupin Limited"
$ws.Cells.Item(22, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(22, 4).Value = "Complete Exit"
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 1.043005
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = -1.043005

# row 23
$ws.Cells.Item(23, 1).Value = "INE296A01032"
$ws.Cells.Item(23, 2).Value = "Bajaj Finance Limited"
$ws.Cells.Item(23, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(23, 4).Value = "Complete Exit"
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 6.442494
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = -6.442494

# row 24
$ws.Cells.Item(24, 1).Value = "INE860A01027"
$ws.Cells.Item(24, 2).Value = "HCL Technologies Limited"
$ws.Cells.Item(24, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(24, 4).Value = "Complete Exit"
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 1.618943
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = -1.618943
$ws.Cells.Item(24, 9).Value = 0

# row 25
$ws.Cells.Item(25, 1).Value = "INE271C01023"
$ws.Cells.Item(25, 2).Value = "DLF Limited"
$ws.Cells.Item(25, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(25, 4).Value = "Complete Exit"
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 5.824721
$ws.Cells.Item(25, 7).Value = 6.25173
$ws.Cells.Item(25, 8).Value = -5.824721
$ws.Cells.Item(25, 9).Value = -6.25173

# row 26
$ws.Cells.Item(26, 1).Value = "INE437A01024"
$ws.Cells.Item(26, 2).Value = "Apollo Hospitals Enterprise Ltd"
$ws.Cells.Item(26, 3).Value = "quant Quantamental Fund"
$ws.Cells.Item(26, 4).Value = "Complete Exit"
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 1.997156
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = -1.997156
$ws.Cells.Item(26, 9).Value = 0

